# Daily TGP (terminal gate pricing) refresh: shift each state's date/price
# rows forward by one effective day, per the updated "before/after" cell map.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value2 = 45961
$ws.Range("D8").Value2 = 167.32
$ws.Range("E8").Value2 = 160.94
$ws.Range("F8").Value2 = 170.94
$ws.Range("G8").Value2 = 161.09

$ws.Range("A9").Value2 = 45961
$ws.Range("D9").Value2 = 167.32
$ws.Range("E9").Value2 = 160.94
$ws.Range("F9").Value2 = 170.94
$ws.Range("G9").Value2 = 161.09

$ws.Range("A10").Value2 = 45961
$ws.Range("D10").Value2 = 169.79
$ws.Range("E10").Value2 = 163.42
$ws.Range("F10").Value2 = 173.42
$ws.Range("G10").Value2 = 163.89

$ws.Range("A11").Value2 = 45960
$ws.Range("D11").Value2 = 165.97
$ws.Range("E11").Value2 = 160.19
$ws.Range("F11").Value2 = 170.19
$ws.Range("G11").Value2 = 160.34

$ws.Range("A12").Value2 = 45960
$ws.Range("D12").Value2 = 165.97
$ws.Range("E12").Value2 = 160.19
$ws.Range("F12").Value2 = 170.19
$ws.Range("G12").Value2 = 160.34

$ws.Range("A13").Value2 = 45960
$ws.Range("D13").Value2 = 168.85
$ws.Range("E13").Value2 = 162.95
$ws.Range("F13").Value2 = 172.95
$ws.Range("G13").Value2 = 163.42

$ws.Range("A17").Value2 = 45961
$ws.Range("D17").Value2 = 173.29
$ws.Range("E17").Value2 = 166.35
$ws.Range("F17").Value2 = 176.35

$ws.Range("A18").Value2 = 45960
$ws.Range("D18").Value2 = 172.33
$ws.Range("E18").Value2 = 165.82
$ws.Range("F18").Value2 = 175.82

$ws.Range("A22").Value2 = 45961
$ws.Range("D22").Value2 = 169.0
$ws.Range("E22").Value2 = 162.44
$ws.Range("F22").Value2 = 172.04
$ws.Range("G22").Value2 = 163.62

$ws.Range("A23").Value2 = 45961
$ws.Range("D23").Value2 = 174.56
$ws.Range("E23").Value2 = 167.15
$ws.Range("F23").Value2 = 177.15

$ws.Range("A24").Value2 = 45961
$ws.Range("D24").Value2 = 174.37
$ws.Range("E24").Value2 = 167.32
$ws.Range("F24").Value2 = 177.32

$ws.Range("A25").Value2 = 45961
$ws.Range("D25").Value2 = 175.2
$ws.Range("E25").Value2 = 166.71
$ws.Range("F25").Value2 = 176.71
$ws.Range("G25").Value2 = 166.54

$ws.Range("A26").Value2 = 45961
$ws.Range("D26").Value2 = 173.94
$ws.Range("E26").Value2 = 168.28
$ws.Range("F26").Value2 = 178.28

$ws.Range("A27").Value2 = 45960
$ws.Range("D27").Value2 = 167.87
$ws.Range("E27").Value2 = 161.97
$ws.Range("F27").Value2 = 171.56
$ws.Range("G27").Value2 = 163.14

$ws.Range("A28").Value2 = 45960
$ws.Range("D28").Value2 = 173.62
$ws.Range("E28").Value2 = 166.68
$ws.Range("F28").Value2 = 176.68

$ws.Range("A29").Value2 = 45960
$ws.Range("D29").Value2 = 173.43
$ws.Range("E29").Value2 = 166.86
$ws.Range("F29").Value2 = 176.86

$ws.Range("A30").Value2 = 45960
$ws.Range("D30").Value2 = 174.26
$ws.Range("E30").Value2 = 166.25
$ws.Range("F30").Value2 = 176.25
$ws.Range("G30").Value2 = 166.08

$ws.Range("A31").Value2 = 45960
$ws.Range("D31").Value2 = 172.99
$ws.Range("E31").Value2 = 167.81
$ws.Range("F31").Value2 = 177.82

$ws.Range("A35").Value2 = 45961
$ws.Range("D35").Value2 = 168.04
$ws.Range("E35").Value2 = 160.63
$ws.Range("F35").Value2 = 169.63

$ws.Range("A36").Value2 = 45960
$ws.Range("D36").Value2 = 167.1
$ws.Range("E36").Value2 = 160.17
$ws.Range("F36").Value2 = 169.17

$ws.Range("A40").Value2 = 45961
$ws.Range("D40").Value2 = 173.64
$ws.Range("E40").Value2 = 165.98
$ws.Range("F40").Value2 = 175.98

$ws.Range("A41").Value2 = 45961
$ws.Range("D41").Value2 = 173.35
$ws.Range("E41").Value2 = 166.4
$ws.Range("F41").Value2 = 176.4

$ws.Range("A42").Value2 = 45960
$ws.Range("D42").Value2 = 172.72
$ws.Range("E42").Value2 = 165.48
$ws.Range("F42").Value2 = 175.48

$ws.Range("A43").Value2 = 45960
$ws.Range("D43").Value2 = 172.43
$ws.Range("E43").Value2 = 165.9
$ws.Range("F43").Value2 = 175.9

$ws.Range("A47").Value2 = 45961
$ws.Range("D47").Value2 = 166.31
$ws.Range("E47").Value2 = 161.08
$ws.Range("F47").Value2 = 171.08

$ws.Range("A48").Value2 = 45961
$ws.Range("D48").Value2 = 166.31
$ws.Range("E48").Value2 = 161.26
$ws.Range("F48").Value2 = 171.26

$ws.Range("A49").Value2 = 45960
$ws.Range("D49").Value2 = 165.14
$ws.Range("E49").Value2 = 160.96
$ws.Range("F49").Value2 = 170.96

$ws.Range("A50").Value2 = 45960
$ws.Range("D50").Value2 = 165.14
$ws.Range("E50").Value2 = 161.14
$ws.Range("F50").Value2 = 171.14

$ws.Range("A54").Value2 = 45961
$ws.Range("D54").Value2 = 183.88
$ws.Range("E54").Value2 = 176.17
$ws.Range("F54").Value2 = 186.17

$ws.Range("A55").Value2 = 45961
$ws.Range("D55").Value2 = 171.53
$ws.Range("E55").Value2 = 173.86
$ws.Range("F55").Value2 = 183.86

$ws.Range("A56").Value2 = 45961
$ws.Range("D56").Value2 = 173.91

$ws.Range("A57").Value2 = 45961
$ws.Range("D57").Value2 = 173.65
$ws.Range("E57").Value2 = 168.13

$ws.Range("A58").Value2 = 45961
$ws.Range("D58").Value2 = 169.56
$ws.Range("E58").Value2 = 164.18
$ws.Range("F58").Value2 = 174.18

$ws.Range("A59").Value2 = 45961
$ws.Range("D59").Value2 = 176.29
$ws.Range("E59").Value2 = 174.46

$ws.Range("A60").Value2 = 45960
$ws.Range("D60").Value2 = 182.95
$ws.Range("E60").Value2 = 175.77
$ws.Range("F60").Value2 = 185.77

$ws.Range("A61").Value2 = 45960
$ws.Range("D61").Value2 = 170.6
$ws.Range("E61").Value2 = 173.26
$ws.Range("F61").Value2 = 183.26

$ws.Range("A62").Value2 = 45960
$ws.Range("D62").Value2 = 172.86

$ws.Range("A63").Value2 = 45960
$ws.Range("D63").Value2 = 172.58
$ws.Range("E63").Value2 = 167.53

$ws.Range("A64").Value2 = 45960
$ws.Range("D64").Value2 = 168.49
$ws.Range("E64").Value2 = 163.58
$ws.Range("F64").Value2 = 173.58

$ws.Range("A65").Value2 = 45960
$ws.Range("D65").Value2 = 175.32
$ws.Range("E65").Value2 = 174.04
